$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = 51938
$ws.Cells.Item(2, 2).Value = "Sra. Lavínia Freitas"
$ws.Cells.Item(2, 3).Value = "Juridico"
$ws.Cells.Item(2, 5).Value2 = 4
$ws.Cells.Item(2, 6).Value2 = 45081
$ws.Cells.Item(2, 7).Value2 = 3907.6

# Row 3
$ws.Cells.Item(3, 1).Value2 = 39741
$ws.Cells.Item(3, 2).Value = "Bruna Rocha"
$ws.Cells.Item(3, 3).Value = "P&D"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 45084
$ws.Cells.Item(3, 7).Value2 = 2075.68

# Row 4
$ws.Cells.Item(4, 1).Value2 = 39649
$ws.Cells.Item(4, 2).Value = "Dra. Ana Clara Costela"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Doenca"
$ws.Cells.Item(4, 6).Value2 = 45085
$ws.Cells.Item(4, 7).Value2 = 4194.79

# Row 5
$ws.Cells.Item(5, 1).Value2 = 4527
$ws.Cells.Item(5, 2).Value = "Maria Liz da Cunha"
$ws.Cells.Item(5, 3).Value = "Recursos Humanos"
$ws.Cells.Item(5, 4).Value = "Consulta medica"
$ws.Cells.Item(5, 5).Value2 = 7
$ws.Cells.Item(5, 6).Value2 = 45082
$ws.Cells.Item(5, 7).Value2 = 7381.91

# Row 6
$ws.Cells.Item(6, 1).Value2 = 76910
$ws.Cells.Item(6, 2).Value = "Dra. Bárbara Martins"
$ws.Cells.Item(6, 3).Value = "Recursos Humanos"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value2 = 4
$ws.Cells.Item(6, 7).Value2 = 9245.84

# Row 7
$ws.Cells.Item(7, 1).Value2 = 17876
$ws.Cells.Item(7, 2).Value = "Carlos Eduardo Araújo"
$ws.Cells.Item(7, 3).Value = "Financeiro"
$ws.Cells.Item(7, 4).Value = "Consulta medica"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 45094
$ws.Cells.Item(7, 7).Value2 = 6640.7

# Row 8
$ws.Cells.Item(8, 1).Value2 = 64927
$ws.Cells.Item(8, 2).Value = "Ravy Gomes"
$ws.Cells.Item(8, 3).Value = "Engenharia"
$ws.Cells.Item(8, 4).Value = "Consulta medica"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 45106
$ws.Cells.Item(8, 7).Value2 = 3310.16

# Row 9
$ws.Cells.Item(9, 1).Value2 = 19507
$ws.Cells.Item(9, 2).Value = "Alana das Neves"
$ws.Cells.Item(9, 3).Value = "Engenharia"
$ws.Cells.Item(9, 5).Value2 = 8
$ws.Cells.Item(9, 6).Value2 = 45099
$ws.Cells.Item(9, 7).Value2 = 4478.38

# Row 10
$ws.Cells.Item(10, 1).Value2 = 38452
$ws.Cells.Item(10, 2).Value = "Kaique Castro"
$ws.Cells.Item(10, 5).Value2 = 4
$ws.Cells.Item(10, 6).Value2 = 45078
$ws.Cells.Item(10, 7).Value2 = 3087.43

# Row 11
$ws.Cells.Item(11, 1).Value2 = 24048
$ws.Cells.Item(11, 2).Value = "Julia Fogaça"
$ws.Cells.Item(11, 3).Value = "P&D"
$ws.Cells.Item(11, 4).Value = "Problemas pessoais"
$ws.Cells.Item(11, 6).Value2 = 45086
$ws.Cells.Item(11, 7).Value2 = 6623.87
